$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 9027
$ws.Range("F4").Value = 326
$ws.Range("F6").Value = 744
$ws.Range("F7").Value = 140
$ws.Range("F10").Value = 921
$ws.Range("F11").Value = 4040
$ws.Range("F12").Value = 321
$ws.Range("F13").Value = 199
$ws.Range("F14").Value = 817
$ws.Range("F17").Value = 511
$ws.Range("F20").Value = 1462
$ws.Range("F21").Value = 1374
$ws.Range("F22").Value = 541
$ws.Range("F26").Value = 394
$ws.Range("F31").Value = 790
$ws.Range("F32").Value = 82
$ws.Range("F34").Value = 120
$ws.Range("F39").Value = 210
$ws.Range("F40").Value = 433
$ws.Range("F42").Value = 34

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 91
$ws.Range("F6").Value = 68

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 224

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 224
$ws.Range("F3").Value = 9027
$ws.Range("F4").Value = 326
$ws.Range("F6").Value = 744
$ws.Range("F7").Value = 140
$ws.Range("F10").Value = 921
$ws.Range("F12").Value = 4040
$ws.Range("F13").Value = 321
$ws.Range("F14").Value = 200
$ws.Range("F16").Value = 91
$ws.Range("F17").Value = 817
$ws.Range("F20").Value = 68
$ws.Range("F22").Value = 511
$ws.Range("F26").Value = 1462
$ws.Range("F27").Value = 1374
$ws.Range("F28").Value = 541
$ws.Range("F33").Value = 394
$ws.Range("F37").Value = 791
$ws.Range("F38").Value = 82
$ws.Range("F40").Value = 120
$ws.Range("F44").Value = 210
$ws.Range("F45").Value = 433
$ws.Range("F47").Value = 34
